# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (with per-fund holdings detail, mirroring
# the existing "2021-Q2"/"2021-Q3" sheets) right before the "总计" (totals)
# sheet, and updates the "总计" sheet with a new top row summarizing the
# 2022-Q1 quarter, pushing the existing 2021-Q3 / 2021-Q2 rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: force a numeric-looking string to be stored as TEXT (preserves
# leading zeros / exact decimal text, matching the source data which stores
# these as text rather than numbers).
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($cell, $val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# ---------------------------------------------------------------------------
# 1. Locate the existing "总计" sheet and insert a new sheet right before it.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Reference sheet that already carries the header / index-column styling we
# want to reuse ("2021-Q3"), so the new sheet's look matches exactly.
$styleSrc = $wb.Worksheets.Item("2021-Q3")

# Copy header-row formatting (bold, centered, bordered) onto B1:H1.
$styleSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# Copy index-column formatting onto A2:A9 (8 data rows).
$styleSrc.Range("A2").Copy()
$q1.Range("A2:A9").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Header row for the new sheet.
# ---------------------------------------------------------------------------
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# 3. Fund holdings detail rows (A=index, B=code, C=name, D=size, E=stock
#    position, F=position ratio, G=held value, H=rank).
#    B/D/E/F/G are stored as TEXT (leading zeros / fixed decimal strings);
#    A and H are numbers.
# ---------------------------------------------------------------------------
$rows = @(
    @(0, "005613", "上投摩根富时发达市场REITs指数QDII人民币份额", "4.84", "91.10", "3.67", "0.1776", 5),
    @(1, "005614", "上投摩根富时发达市场REITs指数QDII美钞",       "4.84", "91.10", "3.67", "0.1776", 5),
    @(2, "005615", "上投摩根富时发达市场REITs指数QDII美汇",       "4.84", "91.10", "3.67", "0.1776", 5),
    @(3, "000179", "广发美国房地产指数QDII-人民币",               "2.37", "92.38", "2.88", "0.0683", 7),
    @(4, "000180", "广发美国房地产指数QDII - 美元",               "2.37", "92.38", "2.88", "0.0683", 7),
    @(5, "160140", "南方道琼斯美国精选REIT指数(QDII-LOF)A",       "1.35", "89.10", "3.16", "0.0427", 6),
    @(6, "070031", "嘉实全球房地产(QDII)",                        "0.60", "95.08", "3.60", "0.0216", 3),
    @(7, "160141", "南方道琼斯美国精选REIT指数(QDII-LOF)C",       "0.44", "89.10", "3.16", "0.0139", 6)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $rowNum = $i + 2

    $q1.Cells.Item($rowNum, 1).Value = $r[0]

    Set-TextValue $q1.Cells.Item($rowNum, 2) $r[1]
    $q1.Cells.Item($rowNum, 3).Value = $r[2]
    Set-TextValue $q1.Cells.Item($rowNum, 4) $r[3]
    Set-TextValue $q1.Cells.Item($rowNum, 5) $r[4]
    Set-TextValue $q1.Cells.Item($rowNum, 6) $r[5]
    Set-TextValue $q1.Cells.Item($rowNum, 7) $r[6]

    $q1.Cells.Item($rowNum, 8).Value = $r[7]
}

# ---------------------------------------------------------------------------
# 4. Update the "总计" sheet: push the two existing data rows down one row
#    and insert the new 2022-Q1 summary row at the top (row 2).
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# Give the about-to-be-created row 4 the same index-column style as the
# existing rows (copy from A3, which already carries it).
$totals.Range("A3").Copy()
$totals.Range("A4").PasteSpecial(-4122)

# Row 4 <= old row 3 ("2021-Q2": 1, 2021-Q2, 7, 0.53), now index 2.
$totals.Cells.Item(4, 1).Value = 2
$totals.Cells.Item(4, 2).Value = "2021-Q2"
$totals.Cells.Item(4, 3).Value = 7
$totals.Cells.Item(4, 4).Value = 0.53

# Row 3 <= old row 2 ("2021-Q3": 0, 2021-Q3, 7, 0.52), now index 1.
$totals.Cells.Item(3, 1).Value = 1
$totals.Cells.Item(3, 2).Value = "2021-Q3"
$totals.Cells.Item(3, 3).Value = 7
$totals.Cells.Item(3, 4).Value = 0.52

# Row 2 <= new "2022-Q1" summary row, index 0.
$totals.Cells.Item(2, 1).Value = 0
$totals.Cells.Item(2, 2).Value = "2022-Q1"
$totals.Cells.Item(2, 3).Value = 8
$totals.Cells.Item(2, 4).Value = 0.75

# ---------------------------------------------------------------------------
# 5. Restore the originally-active sheet/selection (creating/renaming sheets
#    shifts Excel's active tab onto them) so the workbook's active tab stays
#    on "2021-Q2", matching the un-touched original view state.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
